# feat: add 2022-Q4 data
#
# Plan:
#  1. The workbook currently has two sheets: "总计" and "2020-Q4".
#  2. Duplicate "2020-Q4" (creates a 3rd sheet right after it, keeping all
#     of its data/format intact) and rename the duplicate back to
#     "2020-Q4" - this preserves the historical Q4-2020 holdings sheet.
#  3. Clear the original "2020-Q4" sheet and refill it with the new
#     2022-Q4 fund-holding data, then rename it to "2022-Q4".
#  4. Update the "总计" (totals) summary sheet: the former row 2
#     (2020-Q4 totals) moves down to row 3, and a new row 2 is added
#     with the 2022-Q4 totals.

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item(1)
$oldQ4   = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# Step 1: duplicate the existing "2020-Q4" sheet so its data survives
# under its own tab, positioned right after it.
# ---------------------------------------------------------------------
$oldQ4.Copy([System.Reflection.Missing]::Value, $oldQ4)
$preserved2020 = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------
# Step 2: turn the original sheet (still holding the 2020-Q4 data) into
# the new "2022-Q4" sheet - rename it first so the "2020-Q4" name is
# freed up for the duplicate created above.
# ---------------------------------------------------------------------
$newQ4 = $oldQ4
$newQ4.Name = "2022-Q4"
$newQ4.Cells.Clear()

$preserved2020.Name = "2020-Q4"

# Header row - copy the existing header formatting (bold + border, same
# style already used on the "总计" sheet) across the row, then fill in
# the new header text.
$summary.Range("B1").Copy()
$newQ4.Range("B1:H1").PasteSpecial(-4122)

$newQ4.Range("B1").Value = "基金代码"
$newQ4.Range("C1").Value = "基金名称"
$newQ4.Range("D1").Value = "基金规模"
$newQ4.Range("E1").Value = "股票总仓位"
$newQ4.Range("F1").Value = "仓位占比"
$newQ4.Range("G1").Value = "持有市值(亿元)"
$newQ4.Range("H1").Value = "仓位排名"

# Row-index column (A) uses the same formatting as "总计"!A2.
$summary.Range("A2").Copy()
$newQ4.Range("A2:A3").PasteSpecial(-4122)

# Text columns (B:G) must stay text so things like leading zeros in fund
# codes and the "0.63"-style decimal strings are preserved verbatim. The
# "@" number format forces literal-text storage; ClearFormats() afterwards
# drops the now-unneeded number format again (text type is kept).
$newQ4.Range("B2:G3").NumberFormat = "@"

$newQ4.Range("A2").Value = 0
$newQ4.Range("B2").Value = "007506"
$newQ4.Range("C2").Value = "华夏中证AH经济蓝筹股票指数C"
$newQ4.Range("D2").Value = "0.63"
$newQ4.Range("E2").Value = "91.96"
$newQ4.Range("F2").Value = "1.13"
$newQ4.Range("G2").Value = "0.0071"
$newQ4.Range("H2").Value = 9

$newQ4.Range("A3").Value = 1
$newQ4.Range("B3").Value = "007505"
$newQ4.Range("C3").Value = "华夏中证AH经济蓝筹股票指数A"
$newQ4.Range("D3").Value = "0.33"
$newQ4.Range("E3").Value = "91.96"
$newQ4.Range("F3").Value = "1.13"
$newQ4.Range("G3").Value = "0.0037"
$newQ4.Range("H3").Value = 9

$newQ4.Range("B2:G3").ClearFormats()

# ---------------------------------------------------------------------
# Step 3: update the "总计" summary sheet - push the 2020-Q4 totals row
# down to row 3 and put the new 2022-Q4 totals in row 2.
# ---------------------------------------------------------------------
$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2020-Q4"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0.04

$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.01
